$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "43.911.91"
$ws.Cells.Item(2, 5).Value = "  +3.65%  "

$ws.Cells.Item(3, 4).Value = "2.254.67"
$ws.Cells.Item(3, 5).Value = "  +1.15%  "

$ws.Cells.Item(4, 5).Value = "  -0.15%  "

$ws.Cells.Item(5, 4).Value = "230.22"
$ws.Cells.Item(5, 5).Value = "  -0.29%  "

$ws.Cells.Item(6, 4).Value = "0.635"
$ws.Cells.Item(6, 5).Value = "  +2.64%  "

$ws.Cells.Item(7, 4).Value = "62.96"
$ws.Cells.Item(7, 5).Value = "  +4.12%  "

$ws.Cells.Item(8, 5).Value = "  -0.01%  "

$ws.Cells.Item(9, 4).Value = "0.449"
$ws.Cells.Item(9, 5).Value = "  +11.27%  "

$ws.Cells.Item(10, 5).Value = "  +13.23%  "

$ws.Cells.Item(11, 4).Value = "57.23"
$ws.Cells.Item(11, 5).Value = "  -0.38%  "

$ws.Cells.Item(12, 5).Value = "  +2.20%  "

$ws.Cells.Item(13, 4).Value = "25.80"
$ws.Cells.Item(13, 5).Value = "  +16.05%  "

$ws.Cells.Item(14, 4).Value = "2.591.03"
$ws.Cells.Item(14, 5).Value = "  +1.13%  "

$ws.Cells.Item(15, 4).Value = "15.56"
$ws.Cells.Item(15, 5).Value = "  +0.04%  "

$ws.Cells.Item(16, 4).Value = "6.15"
$ws.Cells.Item(16, 5).Value = "  +9.56%  "

$ws.Cells.Item(17, 4).Value = "0.846"
$ws.Cells.Item(17, 5).Value = "  +6.14%  "

$ws.Cells.Item(18, 4).Value = "2.237.74"
$ws.Cells.Item(18, 5).Value = "  -0.02%  "

$ws.Cells.Item(19, 4).Value = "43.832.26"
$ws.Cells.Item(19, 5).Value = "  +3.68%  "

$ws.Cells.Item(20, 4).Value = "0.0000101"
$ws.Cells.Item(20, 5).Value = "  +7.48%  "

$ws.Cells.Item(21, 4).Value = "73.22"
$ws.Cells.Item(21, 5).Value = "  +1.38%  "

$ws.Cells.Item(22, 4).Value = "6.04"
$ws.Cells.Item(22, 5).Value = "  -2.21%  "

$ws.Cells.Item(23, 4).Value = "251.81"
$ws.Cells.Item(23, 5).Value = "  +3.02%  "

$ws.Cells.Item(24, 5).Value = "  +0.23%  "

$ws.Cells.Item(25, 4).Value = "2.44"
$ws.Cells.Item(25, 5).Value = "  +1.71%  "

$ws.Cells.Item(26, 4).Value = "2.34"
$ws.Cells.Item(26, 5).Value = "  -1.59%  "

$ws.Cells.Item(27, 4).Value = "3.30"
$ws.Cells.Item(27, 5).Value = "  +23.93%  "

$ws.Cells.Item(28, 4).Value = "10.01"
$ws.Cells.Item(28, 5).Value = "  +3.27%  "

$ws.Cells.Item(29, 4).Value = "171.79"
$ws.Cells.Item(29, 5).Value = "  +1.50%  "

$ws.Cells.Item(30, 5).Value = "  -1.73%  "

$ws.Cells.Item(31, 4).Value = "20.78"
$ws.Cells.Item(31, 5).Value = "  +2.11%  "

$ws.Cells.Item(32, 5).Value = "  -4.92%  "

$ws.Cells.Item(33, 5).Value = "  +3.04%  "

$ws.Cells.Item(34, 4).Value = "0.0688"
$ws.Cells.Item(34, 5).Value = "  +5.77%  "

$ws.Cells.Item(35, 4).Value = "4.76"
$ws.Cells.Item(35, 5).Value = "  +2.80%  "

$ws.Cells.Item(36, 4).Value = "4.85"
$ws.Cells.Item(36, 5).Value = "  -3.29%  "

$ws.Cells.Item(37, 4).Value = "3.83"
$ws.Cells.Item(37, 5).Value = "  +8.45%  "

$ws.Cells.Item(38, 4).Value = "6.50"
$ws.Cells.Item(38, 5).Value = "  +1.93%  "

$ws.Cells.Item(39, 4).Value = "2.31"
$ws.Cells.Item(39, 5).Value = "  -1.06%  "

$ws.Cells.Item(40, 4).Value = "0.0256"
$ws.Cells.Item(40, 5).Value = "  +2.59%  "

$ws.Cells.Item(41, 5).Value = "  -0.06%  "

$ws.Cells.Item(42, 4).Value = "17.48"
$ws.Cells.Item(42, 5).Value = "  +9.51%  "

$ws.Cells.Item(43, 4).Value = "0.000218"
$ws.Cells.Item(43, 5).Value = "  -5.05%  "

$ws.Cells.Item(44, 4).Value = "8.26"
$ws.Cells.Item(44, 5).Value = "  -3.83%  "

$ws.Cells.Item(45, 4).Value = "0.0971"
$ws.Cells.Item(45, 5).Value = "  +1.16%  "

$ws.Cells.Item(46, 2).Value = "Aave"
$ws.Cells.Item(46, 3).Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Cells.Item(46, 4).Value = "97.54"
$ws.Cells.Item(46, 5).Value = "  +0.74%  "

$ws.Cells.Item(47, 2).Value = "TrustWalletToken"
$ws.Cells.Item(47, 3).Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Cells.Item(47, 4).Value = "1.20"
$ws.Cells.Item(47, 5).Value = "  -0.59%  "

$ws.Cells.Item(48, 2).Value = "FTXToken"
$ws.Cells.Item(48, 3).Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Cells.Item(48, 4).Value = "4.38"
$ws.Cells.Item(48, 5).Value = "  +0.00%  "

$ws.Cells.Item(49, 4).Value = "1.442.73"
$ws.Cells.Item(49, 5).Value = "  -1.13%  "

$ws.Cells.Item(50, 2).Value = "Celestia"
$ws.Cells.Item(50, 3).Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Cells.Item(50, 4).Value = "9.99"
$ws.Cells.Item(50, 5).Value = "  +15.40%  "

$ws.Cells.Item(51, 2).Value = "NEARProtocol"
$ws.Cells.Item(51, 3).Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Cells.Item(51, 4).Value = "2.29"
$ws.Cells.Item(51, 5).Value = "  +2.48%  "
